$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 stays the same except the polite_expressions (C2) value is cleared
$ws.Cells.Item(2, 3).Value = ""

# New row 3 (the newly added annotation row)
$ws.Cells.Item(3, 1).Value = "parisk"
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(3, 3).Value = "nan"
$ws.Cells.Item(3, 4).Value = "APC"
$ws.Cells.Item(3, 5).Value = "THE"
$ws.Cells.Item(3, 6).Value = "c8048836-24fe-4e27-95aa-c7cfb58ac155"
$ws.Cells.Item(3, 7).Value = "rkc_hGb0Z_annotated.xlsx"
$ws.Cells.Item(3, 8).Value = "The structure of the global policies used in the experiments should be mentioned somewhere."
$ws.Cells.Item(3, 9).Value = "Correct"
